$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 43: add marca (B43) and estado "en proceso" (C43)
$ws.Range("B43").Value = "Agustina"
$ws.Range("C43").Value = "en proceso"

# Row 44: change C44 from text "en proceso" to numeric 100% (percentage format)
$ws.Range("C44").Value = 1
$ws.Range("C44").NumberFormat = "0%"

# Widen column A to fit new content, remove autofit best-fit sizing
$ws.Columns.Item(1).ColumnWidth = 76.66666666666667

# Move active selection to C44
$ws.Range("C44").Select()
